# Auto-generated COM-interop edit script for 案件情報.xlsx (sheet "ランサーズ")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Clear existing hyperlinks on the sheet; they will be re-added below with
# the correct target URLs once every row has been rewritten into its new
# (shifted) position.
$ws.Hyperlinks.Delete()

# Row 2: 法人向け生成AIサービス(RAG・議事録機能)の設計・開発を支援エンジニア募集(
$ws.Range("A2").Value = '2025-12-15 18:29:04'
$ws.Range("B2").Value = '法人向け生成AIサービス(RAG・議事録機能)の設計・開発を支援エンジニア募集(AI/バックエンド)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5445159'
$ws.Range("G2").Value = 368
$ws.Range("H2").Value = '🔥AI,Ai ◆開発'

# Row 3: B2B向け生成AIサービス(チャット・RAG)の新規開発プロジェクト推進を支援し
$ws.Range("A3").Value = '2025-12-15 18:29:04'
$ws.Range("B3").Value = 'B2B向け生成AIサービス(チャット・RAG)の新規開発プロジェクト推進を支援してくださるPM募集'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5445154'
$ws.Range("G3").Value = 368
$ws.Range("H3").Value = '🔥AI,Ai ◆開発'

# Row 4: 建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエ
$ws.Range("A4").Value = '2025-12-15 18:29:04'
$ws.Range("B4").Value = '建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5434128'
$ws.Range("G4").Value = 368
$ws.Range("H4").Value = '🔥AI,Ai ◆開発'

# Row 5: 企業のMicrosoft Copilot導入・活用支援AIコンサルタント募集(研
$ws.Range("A5").Value = '2025-12-15 18:29:04'
$ws.Range("B5").Value = '企業のMicrosoft Copilot導入・活用支援AIコンサルタント募集(研修講師・メンター)'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5434363'
$ws.Range("G5").Value = 348
$ws.Range("H5").Value = '🔥AI,Ai ◆コンサル'

# Row 6: 【AIシステム構築】次のテストに向けた宿題自動出題システム
$ws.Range("A6").Value = '2025-12-15 18:29:04'
$ws.Range("B6").Value = '【AIシステム構築】次のテストに向けた宿題自動出題システム'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5453785'
$ws.Range("G6").Value = 318
$ws.Range("H6").Value = '🔥AI,Ai'

# Row 7: 大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(
$ws.Range("A7").Value = '2025-12-15 18:29:04'
$ws.Range("B7").Value = '大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5427956'
$ws.Range("G7").Value = 310
$ws.Range("H7").Value = '🔥AI,Ai'

# Row 8: 【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像
$ws.Range("A8").Value = '2025-12-15 18:29:04'
$ws.Range("B8").Value = '【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像認識/動画解析)'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5439158'
$ws.Range("G8").Value = 303
$ws.Range("H8").Value = '🔥AI,Ai'

# Row 9: AIオートメーションエンジニア
$ws.Range("A9").Value = '2025-12-15 18:29:04'
$ws.Range("B9").Value = 'AIオートメーションエンジニア'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5453810'
$ws.Range("G9").Value = 298
$ws.Range("H9").Value = '🔥AI,Ai'

# Row 10: 海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動
$ws.Range("A10").Value = '2025-12-15 18:29:04'
$ws.Range("B10").Value = '海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5251319'
$ws.Range("G10").Value = 135
$ws.Range("H10").Value = '◆ツール,スクレイピング ◇サイト'

# Row 11: 【Unity/XRエンジニア募集】製造業DX支援!既存システムと連携するXRアプ
$ws.Range("A11").Value = '2025-12-15 18:29:04'
$ws.Range("B11").Value = '【Unity/XRエンジニア募集】製造業DX支援!既存システムと連携するXRアプリ開発'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5454210'
$ws.Range("G11").Value = 108
$ws.Range("H11").Value = '◆開発 ◇アプリ'

# Row 12: Javaプログラミング研修の演習サポート講師業務【経験不問】(再掲)
$ws.Range("A12").Value = '2025-12-15 18:29:04'
$ws.Range("B12").Value = 'Javaプログラミング研修の演習サポート講師業務【経験不問】(再掲)'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5453723'
$ws.Range("G12").Value = 85
$ws.Range("H12").Value = '★Java'

# Row 13: GoogleCloudを利用したアジャイル開発共通基盤のSREエンジニアの募集
$ws.Range("A13").Value = '2025-12-15 18:29:04'
$ws.Range("B13").Value = 'GoogleCloudを利用したアジャイル開発共通基盤のSREエンジニアの募集'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5453768'
$ws.Range("G13").Value = 75
$ws.Range("H13").Value = '◆開発'

# Row 14: Base無在庫ツール作成 経験者のみ募集 実績提示をお願いします
$ws.Range("A14").Value = '2025-12-15 18:29:04'
$ws.Range("B14").Value = 'Base無在庫ツール作成 経験者のみ募集 実績提示をお願いします'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5453611'
$ws.Range("G14").Value = 73
$ws.Range("H14").Value = '◆ツール'

# Row 15: クラウド(AWS/Azure) 運用管理 研修の演習サポート講師業務【経験不問】
$ws.Range("A15").Value = '2025-12-15 18:29:04'
$ws.Range("B15").Value = 'クラウド(AWS/Azure) 運用管理 研修の演習サポート講師業務【経験不問】(再掲)'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5453718'
$ws.Range("G15").Value = 38
$ws.Range("H15").Value = '◇管理'

# Row 16: 【介護事業所向け】グーグルワークスペース社内システム構築依頼
$ws.Range("A16").Value = '2025-12-15 18:29:04'
$ws.Range("B16").Value = '【介護事業所向け】グーグルワークスペース社内システム構築依頼'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5453868'
$ws.Range("G16").Value = 40
$ws.Range("H16").ClearContents()

# Row 17: 【急募】ネイティブjs案件 長期でご依頼できるパートナー大募集
$ws.Range("A17").Value = '2025-12-15 18:29:04'
$ws.Range("B17").Value = '【急募】ネイティブjs案件 長期でご依頼できるパートナー大募集'
$ws.Range("C17").Value = 'システム開発'
$ws.Range("D17").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E17").Value = '期限情報なし'
$ws.Range("F17").Value = 'https://www.lancers.jp/work/detail/5454504'
$ws.Range("G17").Value = 25
$ws.Range("H17").ClearContents()

# Row 18: 【急募】ネイティブjsのスペシャリスト募集!
$ws.Range("A18").Value = '2025-12-15 18:29:04'
$ws.Range("B18").Value = '【急募】ネイティブjsのスペシャリスト募集!'
$ws.Range("C18").Value = 'システム開発'
$ws.Range("D18").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E18").Value = '期限情報なし'
$ws.Range("F18").Value = 'https://www.lancers.jp/work/detail/5454495'
$ws.Range("G18").Value = 18
$ws.Range("H18").ClearContents()

# Row 19: 【急募】ホームページとLPの改善をお手伝いします!
$ws.Range("A19").Value = '2025-12-15 18:29:04'
$ws.Range("B19").Value = '【急募】ホームページとLPの改善をお手伝いします!'
$ws.Range("C19").Value = 'システム開発'
$ws.Range("D19").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E19").Value = '期限情報なし'
$ws.Range("F19").Value = 'https://www.lancers.jp/work/detail/5453763'
$ws.Range("G19").Value = 18
$ws.Range("H19").ClearContents()

# Row 20: サンプルER図の作成(研修用にER図はどんなものかをサンプルでみせるもの)
$ws.Range("A20").Value = '2025-12-15 18:29:04'
$ws.Range("B20").Value = 'サンプルER図の作成(研修用にER図はどんなものかをサンプルでみせるもの)'
$ws.Range("C20").Value = 'システム開発'
$ws.Range("D20").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E20").Value = '期限情報なし'
$ws.Range("F20").Value = 'https://www.lancers.jp/work/detail/5454604'
$ws.Range("G20").Value = 10
$ws.Range("H20").ClearContents()

# Re-create hyperlinks for F2:F20 pointing at the (new) target URLs.
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5445159') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5445154') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5434128') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5434363') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5453785') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5427956') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5439158') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5453810') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5251319') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5454210') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5453723') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5453768') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5453611') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5453718') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5453868') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5454504') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F18"), 'https://www.lancers.jp/work/detail/5454495') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F19"), 'https://www.lancers.jp/work/detail/5453763') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F20"), 'https://www.lancers.jp/work/detail/5454604') | Out-Null

# Column B widened from 51 to 52 (raw OOXML width). Excels ColumnWidth
# property is offset from the raw column width by the font padding constant
# (~0.83 for the workbook default font), so 52 - 0.83 = 51.17.
$ws.Columns.Item(2).ColumnWidth = 51.17

$wb.Save()
